$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet2"

$ws.Range("A1").Value = "EmpID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Title"

$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "Tech Global"
$ws.Range("C2").Value = "DevOps"

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "Ulan"
$ws.Range("C3").Value = "Developer"

$ws.Range("A4").Value = 103
$ws.Range("B4").Value = "Abe"
$ws.Range("C4").Value = "Instructor"
